$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 428 (pushes existing rows 428-473 down to 429-474,
# matching the diff's dimension change from A1:R473 to A1:R474).
$ws.Rows(428).Insert()

# Populate the newly inserted row 428 with the new weekly data point.
$ws.Range("A428").Value = 4
$ws.Range("B428").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C428").Value = "Los Lagos"
$ws.Range("D428").Value = 44946
$ws.Range("E428").Value = 10
$ws.Range("F428").Value = 100112008
$ws.Range("G428").Value = "Coliflor"
$ws.Range("H428").Value = "Sin especificar"
$ws.Range("I428").Value = "Primera"
$ws.Range("J428").Value = 1200
$ws.Range("K428").Value = 1600
$ws.Range("L428").Value = 1600
$ws.Range("M428").Value = 1600
$ws.Range("N428").Value = "$/unidad"
$ws.Range("O428").Value = "Región Metropolitana"
$ws.Range("P428").Value = 1600
$ws.Range("Q428").Value = 1
$ws.Range("R428").Value = "Hortaliza"
